$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.350.79"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.871.45"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.83"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4690"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2846"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06546"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.47"
$ws.Range("E10").Value = "  +8.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07891"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.29"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "1.872.97"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.118"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6768"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.51"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "30.339.01"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +2.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.479"
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("D21").Value = "2.114.24"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007321"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.153"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.43"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.179"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.23"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.933"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.387"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09723"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.405"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04719"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("E35").Value = "  +5.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7068"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.726"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.335"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.539"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.95"
$ws.Range("E41").Value = "  +5.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.950"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8515"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4196"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.87"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.233"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.234"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "941.86"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05643"
$ws.Range("E51").Value = "  -0.49%  "
